$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.970.67'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.228.77'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.30'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +7.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.622'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.35'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.79%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.03'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +13.71%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.16'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.16'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +6.35%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.561.38'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.866'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.227.14'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.843.85'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0968'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.22'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.93'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.40'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.20'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +14.49%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.53'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.69%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.24'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.24%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.77'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.83%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.57'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0733'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.74'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.24'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +21.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.96'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +8.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0305'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +13.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.29'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.99'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '67.06'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.05'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +18.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.94'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.203'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +6.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.77'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.66'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.21%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +6.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.19'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.24%  '
